# Fruta / hortaliza, semanal
# Insert a new weekly data row above the current row 54, shifting the
# existing rows 54-60 down to 55-61, then populate the newly inserted
# row 54 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 54; existing rows 54:60 shift down to 55:61.
$ws.Rows.Item(54).Insert()

# Populate the new row 54 with the new weekly observation.
$ws.Cells.Item(54, 1).Value = 10
$ws.Cells.Item(54, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(54, 3).Value = "La Araucanía"
$ws.Cells.Item(54, 4).Value = 44522
$ws.Cells.Item(54, 5).Value = 9
$ws.Cells.Item(54, 6).Value = "Fruta"
$ws.Cells.Item(54, 7).Value = 100101
$ws.Cells.Item(54, 8).Value = "Berries"
$ws.Cells.Item(54, 9).Value = 100101001
$ws.Cells.Item(54, 10).Value = "Arándano (blue)"
$ws.Cells.Item(54, 11).Value = "Sin especificar"
$ws.Cells.Item(54, 12).Value = "Primera"
$ws.Cells.Item(54, 13).Value = 200
$ws.Cells.Item(54, 14).Value = 3200
$ws.Cells.Item(54, 15).Value = 3200
$ws.Cells.Item(54, 16).Value = 3200
$ws.Cells.Item(54, 17).Value = "$/kilo"
$ws.Cells.Item(54, 18).Value = "Región del Maule"
$ws.Cells.Item(54, 19).Value = 3200
$ws.Cells.Item(54, 20).Value = 1
